# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Swap the "Periodo Mora" / "Valor Mora" figures between row 16 and row 17
# so that period 1711 (10824) now appears before period 1712 (29520).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: Periodo Mora / Valor Mora
$ws.Range("E16").Value = "1711"
$ws.Range("F16").Value = 10824

# Row 17: Periodo Mora / Valor Mora
$ws.Range("E17").Value = "1712"
$ws.Range("F17").Value = 29520
